# MonsterAll2015.xlsx edit script
# - Renames Sheet1 -> Teams
# - Adds a new "Side" sheet with the tiered side-bet payout calculator
# - Teams stays first, Side is added after it and becomes the active sheet

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet, add the new one right after it -------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Teams"

$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Side"

# --- Column widths (bestFit, matches the authored layout) ------------------
$ws2.Range("A1").ColumnWidth = 10.42578125
$ws2.Range("B1").ColumnWidth = 14.5703125
$ws2.Range("C1").ColumnWidth = 10.5703125
$ws2.Range("E1:G1").ColumnWidth = 10.85546875
$ws2.Range("H1").ColumnWidth = 13.28515625

# --- Enter the labels/values in the same order the original author typed
#     them, so new shared-string entries land at the same indices.
$ws2.Range("D1").Value = "pay spots"
$ws2.Range("C2").Value = "total purse"
$ws2.Range("A3").Value = "per team"
$ws2.Range("H4").Value = "pay per team"
$ws2.Range("H3").Value = "pay per tier"
$ws2.Range("H2").Value = "pay % per tier"
$ws2.Range("B1").Value = "tiered side bet"
$ws2.Range("A2").Value = "# teams"

# --- Formulas / numbers --------------------------------------------------
$ws2.Range("E1").Formula = "=(CONCATENATE(""pay spots "",`$D`$2/3))"
$ws2.Range("F1:G1").Formula = "=(CONCATENATE(""pay spots "",`$D`$2/3))"

$ws2.Range("B2").Value = 100
$ws2.Range("D2").Formula = "=(3*(FLOOR.MATH(B2/50)-1))+9"
$ws2.Range("E2").Value = 0.6
$ws2.Range("F2").Value = 0.3
$ws2.Range("G2").Value = 0.1

$ws2.Range("B3").Value = 40
$ws2.Range("C3").Formula = "=B3*B2"
$ws2.Range("E3").Formula = "=C3*E2"
$ws2.Range("F3").Formula = "=C3*F2"
$ws2.Range("G3").Formula = "=C3*G2"

$ws2.Range("E4").Formula = "=E3/(`$D`$2/3)"
$ws2.Range("F4:G4").Formula = "=F3/(`$D`$2/3)"

# --- Number formats (build currency-ish #,##0.00 used across the calc) -----
$ws2.Range("B3:G3").NumberFormat = "#,##0.00"
$ws2.Range("B4:G4").NumberFormat = "#,##0.00"

# --- Protection: lock the label/output cells, leave the two inputs
#     (# teams, per team) unlocked so players can edit them -----------------
$ws2.Range("C1:H1").Locked = $true
$ws2.Range("C2:H2").Locked = $true
$ws2.Range("H3:H4").Locked = $true
$ws2.Range("C3:G3").Locked = $true
$ws2.Range("C4:G4").Locked = $true

$ws2.Range("B2").Locked = $false
$ws2.Range("B3").Locked = $false

# --- Sheet view / selection -------------------------------------------
$ws2.Range("E4").Select()

Write-Host "done"
